$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id = 1)
$ws.Range("B2").Value = 76.95
$ws.Range("C2").Value = "The applicant possesses strong skills in Python, JavaScript, ReactJS, and Django, which align well with the requirements for the SDE Intern position. The applicant's experience in developing web and mobile applications using ReactJS and Django, and utilizing PyTorch, Sklearn, and Python for data analysis demonstrates a solid technical foundation. However, the applicant lacks direct experience with MongoDB and NodeJS, which are key requirements for the role. While the applicant has demonstrated proficiency in similar technologies, additional training or experience in MongoDB and NodeJS may be beneficial for ensuring a seamless fit into the role. Overall, the applicant's technical expertise and project experience make them a strong candidate for the position, with the potential for further development in MongoDB and NodeJS."
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "Based on the applicant's responses, it's evident that they exhibit strong adaptability, problem-solving skills, and a willingness to learn and improve. Their expressed interest in learning from the Japanese work culture, seeking feedback, and actively participating in relevant workshops and online courses demonstrates a proactive approach to personal and professional growth. The applicant's ability to overcome setbacks and their curiosity-driven mindset further reinforce their adaptability and problem-solving acumen. Additionally, their articulated vision for contributing to the company's goals, expectations from the company, and desire to immerse themselves in Japanese culture reflect a positive attitude and alignment with the company's values. The applicant's responses indicate a high level of motivation, adaptability, and a sincere interest in working in Japan."

# Row 3 (id = 2)
$ws.Range("B3").Value = 89.06999999999999
$ws.Range("C3").Value = "The applicant has worked on projects that showcase his expertise in web and mobile application development. While the technologies used in the projects are not an exact match with the required skills of MongoDB, ExpressJS, ReactJS, and NodeJS, the applicant's skills in Flutter, Dart, WebRTC, Socket.IO, and HTML demonstrate an aptitude for web and mobile development. However, there is a gap in the required skills as the applicant lacks proficiency in MongoDB and JavaScript, which are essential for the role. Despite this, the applicant's diverse skill set and experience in developing functional applications display potential for growth and adaptability to the technologies used by the company."
$ws.Range("E3").Value = "The applicant has expressed a deep interest in working in Japan, citing the country's leading position in AI and expressing fascination with the work culture. The applicant has also demonstrated openness to adapting to new cultures and learning a new language while working in Japan. Furthermore, the applicant has emphasized the desire to work with a team and the willingness to learn new skills, aligning with the company's soft skills requirements. The applicant's proactive approach to addressing weaknesses and ability to bridge technical and creative aspects in projects reflect an adaptable and team-oriented nature. Overall, the applicant's responses indicate a strong willingness to work in Japan and possess the necessary soft skills for the job."

# Row 4 (id = 3)
$ws.Range("B4").Value = 78.72
$ws.Range("C4").Value = "The applicant possesses strong skills in Python, PyTorch, Tensorflow, and ReactJS, demonstrating proficiency in web and mobile application development. While the applicant's projects showcase extensive experience in computer vision, image processing, and backend development, there is a lack of direct experience with MongoDB, ExpressJS, and NodeJS, which are essential for the role. However, the applicant's versatility and experience in developing web applications using ReactJS could contribute to the role's responsibilities. Despite the gaps in experience with specific technologies, the applicant's strong foundation in relevant areas positions them as a promising candidate for the SDE Intern position."
$ws.Range("E4").Value = "The applicant has exhibited a keen interest in being part of Japan's work culture, showcasing adaptability and a willingness to learn a new language while working in Japan. The applicant's responses convey a strong sense of commitment to pursuing a stable career in the field of AI/ML and backend development, aligning with the company's requirements. Furthermore, the applicant's acknowledgement of strengths in teamwork and openness about areas for improvement demonstrates a growth mindset and a potential for future development. Overall, the applicant's responses indicate a positive attitude and willingness to embrace the opportunities presented in Japan."

# Row 5 (id = 4)
$ws.Range("B5").Value = 70.89
$ws.Range("C5").Value = "The applicant has worked on various web development projects such as a Video Conferencing Project, TEDxIITGuwahati Website, and Alcheringa Pass Portal, showcasing strong skills in JavaScript, HTML, and CSS. However, the job requires expertise in MongoDB, ReactJS, and NodeJS, which are not fully exhibited in the applicant's projects. While the applicant's experience in web development is commendable, there is a lack of direct experience with the specific technologies mentioned in the job description. This indicates potential gaps in the applicant's skillset with regards to the technologies required for the job."
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "The applicant has demonstrated an eagerness to work in Japan, expressing a genuine interest in the country's work culture and a willingness to adapt to a new environment. The applicant's openness to working in diverse locations demonstrates adaptability. Furthermore, the applicant's acknowledgment of strengths in teamwork aligns with the company's soft skill requirements. However, the applicant's admission of struggles with acquiring new skills and preference for working in a team may indicate a need for further development in individual problem-solving and time management skills, resulting in a score of 3 out of 5."

$wb.Save()
